$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing values on row 45 (C45 and D45)
$ws.Range("C45").Value = 233.68
$ws.Range("D45").Value = 1589.03

# Add a new row 48 for period "III-2021"
$ws.Range("A48").Value = "III-2021"
$ws.Range("B48").Value = 8345.24
$ws.Range("C48").Value = 254.08
$ws.Range("D48").Value = 1791.54
$ws.Range("E48").Value = 6028.97
$ws.Range("F48").Value = 212.04
$ws.Range("G48").Value = 58.61
